# Add a "random forest" row plus two new metric columns (cross val mean,
# standard dev) to the scoring sheet, and restyle the numeric cells to a
# uniform 0.0000 number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (D1:E1), bold like the rest of row 1 ------------------
$ws.Range("D1").Value = "cross val mean"
$ws.Range("E1").Value = "standard dev"
$ws.Range("B1:E1").Font.Bold = $true

# --- New data: cross val mean / standard dev for the existing 3 methods -----
$ws.Range("D2").Value = -0.165418449910799
$ws.Range("E2").Value = 0.0202108777533173

$ws.Range("D3").Value = -0.22805646190949
$ws.Range("E3").Value = 0.021937123415589

$ws.Range("D4").Value = -0.169704381240851
$ws.Range("E4").Value = 0.0229859813789266

# --- New row: random forest ---------------------------------------------------
$ws.Range("A5").Value = "random forest"
$ws.Range("B5").Value = 0.214097248398039
$ws.Range("C5").Value = 0.462706438682281
$ws.Range("D5").Value = -0.16325451012216
$ws.Range("E5").Value = 0.0233041525096976

# --- Number formats: everything numeric now uses 0.0000 ---------------------
$ws.Range("B1:E5").NumberFormat = "0.0000"

# --- Column widths for the new columns ---------------------------------------
$ws.Columns.Item(4).ColumnWidth = 15.5
$ws.Columns.Item(5).ColumnWidth = 12.5

# --- Selection moves to E5, matching the saved workbook state ---------------
$ws.Range("E5").Select() | Out-Null
